# Daily "Förändrad" (changed) timestamp roll-forward + drop of the oldest
# tracked case row, mirroring the automatic update performed by the
# upstream generator (openpyxl) each day:
#   - Column C ("Förändrad") advances by one day (45538 -> 45539) for every
#     data row that survives the update.
#   - The last row of data (row 29, "A 36712-2024") is removed - it was the
#     newest entry and drops out of the rolling window.
#   - After the removal, the new last row (28) loses its explicit row
#     height / customHeight flag, same as the previous last row always did
#     before another row was appended after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date for every existing data row (2-29) by one day
# before the last row is dropped.
$ws.Range("C2:C29").Value = 45539

# Drop the last data row (row 29) entirely - content, formatting, everything.
$ws.Rows(29).Delete()

# The new last row (28) reverts to the sheet's standard (non-custom) row
# height, matching the pre-update pattern where only the newest row lacked
# an explicit height.
$ws.Rows(28).AutoFit()
